{"js": "const replacements = [\n  [\"2024-09-06 Friday\", \"2024-09-07 Saturday\"],\n  [\"90\u00f72=45, 0\", \"79\u00f77=11, 2\"],\n  [\"83\u00f78=10, 3\", \"33\u00f74=8, 1\"],\n  [\"91\u00f72=45, 1\", \"74\u00f79=8, 2\"],\n  [\"85\u00f76=14, 1\", \"34\u00f78=4, 2\"],\n  [\"35\u00f73=11, 2\", \"32\u00f79=3, 5\"],\n  [\"62\u00f75=12, 2\", \"34\u00f77=4, 6\"],\n  [\"52\u00f74=13, 0\", \"58\u00f72=29, 0\"],\n  [\"19\u00f74=4, 3\", \"60\u00f78=7, 4\"],\n  [\"13\u00f77=1, 6\", \"63\u00f75=12, 3\"],\n  [\"37\u00f73=12, 1\", \"22\u00f77=3, 1\"],\n  [\"50\u00f77=7, 1\", \"73\u00f77=10, 3\"],\n  [\"16\u00f79=1, 7\", \"98\u00f77=14, 0\"],\n  [\"89\u00f72=44, 1\", \"79\u00f73=26, 1\"],\n  [\"97\u00f79=10, 7\", \"65\u00f72=32, 1\"],\n  [\"43\u00f75=8, 3\", \"46\u00f74=11, 2\"],\n  [\"12\u00f72=6, 0\", \"39\u00f72=19, 1\"],\n  [\"61\u00f77=8, 5\", \"95\u00f73=31, 2\"],\n  [\"51\u00f77=7, 2\", \"11\u00f79=1, 2\"],\n  [\"64\u00f73=21, 1\", \"35\u00f78=4, 3\"],\n  [\"86\u00f72=43, 0\", \"42\u00f76=7, 0\"],\n  [\"49\u00f76=8, 1\", \"11\u00f74=2, 3\"],\n  [\"97\u00f78=12, 1\", \"30\u00f79=3, 3\"],\n  [\"47\u00f72=23, 1\", \"26\u00f78=3, 2\"],\n  [\"67\u00f77=9, 4\", \"80\u00f72=40, 0\"],\n  [\"98\u00f72=49, 0\", \"23\u00f76=3, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2024-09-06 Friday\", \"2024-09-07 Saturday\")\n    ,@(\"90\u00f72=45, 0\", \"79\u00f77=11, 2\")\n    ,@(\"83\u00f78=10, 3\", \"33\u00f74=8, 1\")\n    ,@(\"91\u00f72=45, 1\", \"74\u00f79=8, 2\")\n    ,@(\"85\u00f76=14, 1\", \"34\u00f78=4, 2\")\n    ,@(\"35\u00f73=11, 2\", \"32\u00f79=3, 5\")\n    ,@(\"62\u00f75=12, 2\", \"34\u00f77=4, 6\")\n    ,@(\"52\u00f74=13, 0\", \"58\u00f72=29, 0\")\n    ,@(\"19\u00f74=4, 3\", \"60\u00f78=7, 4\")\n    ,@(\"13\u00f77=1, 6\", \"63\u00f75=12, 3\")\n    ,@(\"37\u00f73=12, 1\", \"22\u00f77=3, 1\")\n    ,@(\"50\u00f77=7, 1\", \"73\u00f77=10, 3\")\n    ,@(\"16\u00f79=1, 7\", \"98\u00f77=14, 0\")\n    ,@(\"89\u00f72=44, 1\", \"79\u00f73=26, 1\")\n    ,@(\"97\u00f79=10, 7\", \"65\u00f72=32, 1\")\n    ,@(\"43\u00f75=8, 3\", \"46\u00f74=11, 2\")\n    ,@(\"12\u00f72=6, 0\", \"39\u00f72=19, 1\")\n    ,@(\"61\u00f77=8, 5\", \"95\u00f73=31, 2\")\n    ,@(\"51\u00f77=7, 2\", \"11\u00f79=1, 2\")\n    ,@(\"64\u00f73=21, 1\", \"35\u00f78=4, 3\")\n    ,@(\"86\u00f72=43, 0\", \"42\u00f76=7, 0\")\n    ,@(\"49\u00f76=8, 1\", \"11\u00f74=2, 3\")\n    ,@(\"97\u00f78=12, 1\", \"30\u00f79=3, 3\")\n    ,@(\"47\u00f72=23, 1\", \"26\u00f78=3, 2\")\n    ,@(\"67\u00f77=9, 4\", \"80\u00f72=40, 0\")\n    ,@(\"98\u00f72=49, 0\", \"23\u00f76=3, 5\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $found = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
